# Update a question text in the itemBank sheet.
# "What is a group of crows called?" / "Murder"  -->
# 'Is it "I have less photos than last year" or "I have fewer photos than last year"' / 'fewer'

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("itemBank")

$ws.Range("C26").Value = 'Is it "I have less photos than last year" or "I have fewer photos than last year"'
$ws.Range("D26").Value = "fewer"

# Reflect the author's on-save view state (scrolled down, new selection).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("C25").Select()
